# Trade #35 closed at 2026-02-16 21:28:40 - momentum DOWN +0.000%
# Appends a new "OPEN" trade row (row 7) to the "momentum" sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("momentum")

$ws.Range("A7").Value = 35

# Force the date-like text to stay a literal string instead of being
# auto-converted to a date serial by the Value setter's type inference,
# then drop back to the default style so no stray formatting is left on
# the cell (matches how every other "Date" cell in this column is stored).
$ws.Range("B7").NumberFormat = "@"
$ws.Range("B7").Value = "2026-02-16"
$ws.Range("B7").Style = "Normal"

$ws.Range("C7").Value = "21:28:40"
$ws.Range("D7").Value = "momentum"
$ws.Range("E7").Value = "DOWN"
$ws.Range("F7").Value = 68718.965
# G7 (Exit Price) intentionally left blank - trade is still OPEN, no exit yet.
$ws.Range("H7").Value = "OPEN"
$ws.Range("I7").Value = 0
$ws.Range("J7").Value = 0
$ws.Range("K7").Value = 0.9
$ws.Range("L7").Value = "Downward momentum: -0.299% over 10 samples"
# M7 (Exit Reason) intentionally left blank - trade is still OPEN, no exit yet.
$ws.Range("N7").Value = 0
